$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DataCombined")

$ws.Range("B2").Value = "simulated"
$ws.Range("C2").Value = "Aciclovir simulated"
$ws.Range("D2").Value = "TestScenario"

$ws.Range("D2").Select()
